$wb = $excel.ActiveWorkbook

# --- Sheet 1: ROW50-FE-LIFTER -> append row 53 ---
$ws1 = $wb.Worksheets.Item(1)
$r = 53
$ws1.Cells.Item($r, 1).Value = 45750.69791148148
$ws1.Cells.Item($r, 1).NumberFormat = $ws1.Cells.Item($r - 1, 1).NumberFormat
$ws1.Cells.Item($r, 2).Value = "0x01,0x90"
$ws1.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws1.Cells.Item($r, 4).Value = "0x01,0x62"
$ws1.Cells.Item($r, 5).Value = "0xe"
$ws1.Cells.Item($r, 6).Value = 400
$ws1.Cells.Item($r, 7).Value = [double]"5.68631262647114e+23"
$ws1.Cells.Item($r, 8).Value = 354
$ws1.Cells.Item($r, 9).Value = 14

# --- Sheet 2: ROW50-MID-LIFTER -> append row 55 ---
$ws2 = $wb.Worksheets.Item(2)
$r = 55
$ws2.Cells.Item($r, 1).Value = 45750.66854166667
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item($r - 1, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = "0x01,0x90 "
$ws2.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws2.Cells.Item($r, 4).Value = "0x01,0x66"
$ws2.Cells.Item($r, 5).Value = "0x19"
$ws2.Cells.Item($r, 6).Value = 400
# Force this long digit-string to be stored as text (matches the rest of
# column G on this sheet) instead of being auto-parsed as a number: set a
# text number-format just for the assignment, then drop back to the
# worksheet's default style so the cell ends up unstyled, same as its peers.
$ws2.Cells.Item($r, 7).NumberFormat = "@"
$ws2.Cells.Item($r, 7).Value = "568631262647113771663628"
$ws2.Cells.Item($r, 7).Style = "Normal"
$ws2.Cells.Item($r, 8).Value = 358
$ws2.Cells.Item($r, 9).Value = 25

# --- Sheet 3: ROW11-FE-LIFTER -> append row 53 ---
$ws3 = $wb.Worksheets.Item(3)
$r = 53
$ws3.Cells.Item($r, 1).Value = 45750.72861877315
$ws3.Cells.Item($r, 1).NumberFormat = $ws3.Cells.Item($r - 1, 1).NumberFormat
$ws3.Cells.Item($r, 2).Value = "0x01,0x90"
$ws3.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws3.Cells.Item($r, 4).Value = "0x01,0x62"
$ws3.Cells.Item($r, 5).Value = "0x14"
$ws3.Cells.Item($r, 6).Value = 400
$ws3.Cells.Item($r, 7).Value = [double]"5.68631262647114e+23"
$ws3.Cells.Item($r, 8).Value = 354
$ws3.Cells.Item($r, 9).Value = 20

# --- Sheet 4: ROW11-MID-LIFTER -> append row 53 ---
$ws4 = $wb.Worksheets.Item(4)
$r = 53
$ws4.Cells.Item($r, 1).Value = 45750.8616753125
$ws4.Cells.Item($r, 1).NumberFormat = $ws4.Cells.Item($r - 1, 1).NumberFormat
$ws4.Cells.Item($r, 2).Value = "0x01,0x90"
$ws4.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws4.Cells.Item($r, 4).Value = "0x01,0x6a"
$ws4.Cells.Item($r, 5).Value = "0x19"
$ws4.Cells.Item($r, 6).Value = 400
$ws4.Cells.Item($r, 7).Value = [double]"5.68631262647114e+23"
$ws4.Cells.Item($r, 8).Value = 362
$ws4.Cells.Item($r, 9).Value = 25
